# Auto-generated Excel COM-interop edit script.
# Updates per-leve market-price snapshot columns (H: currentAveragePrice,
# I: currentAveragePriceNQ, J: currentAveragePriceHQ, K: LevePriceNQ,
# L: LevePriceHQ, M: LeveProfitNQ, N: LeveProfitHQ) across all 8 sheets,
# matching the scheduled market-data refresh described in the commit.

$wb = $excel.ActiveWorkbook


# ==== Sheet: ALC ====
$ws = $wb.Worksheets.Item("ALC")

# Row 31
$ws.Range("H31").Value = 144
$ws.Range("I31").Value = 144
$ws.Range("K31").Value = 432
$ws.Range("M31").Value = -202
# Row 43
$ws.Range("H43").Value = 2675.2222
$ws.Range("J43").Value = 2179.3333
$ws.Range("L43").Value = 2179.3333
$ws.Range("N43").Value = -2317.3333
# Row 98
$ws.Range("H98").Value = 1117.8
$ws.Range("I98").Value = 647.25
$ws.Range("K98").Value = 647.25
$ws.Range("M98").Value = 850.75
# Row 100
$ws.Range("H100").Value = 2750
$ws.Range("I100").Value = 2333.3333
$ws.Range("K100").Value = 2333.3333
$ws.Range("M100").Value = -1792.3333
# Row 115
$ws.Range("H115").Value = 615.3333
$ws.Range("I115").Value = 615.3333
$ws.Range("K115").Value = 1845.9999
$ws.Range("M115").Value = -278.9999
# Row 122
$ws.Range("H122").Value = 1117.8
$ws.Range("I122").Value = 647.25
$ws.Range("K122").Value = 1941.75
$ws.Range("M122").Value = 508.25
# Row 132
$ws.Range("H132").Value = 2809.5
$ws.Range("I132").Value = 1782.4286
$ws.Range("K132").Value = 5347.2858
$ws.Range("M132").Value = -2817.2858

# ==== Sheet: ARM ====
$ws = $wb.Worksheets.Item("ARM")

# Row 2
$ws.Range("H2").Value = 2672
$ws.Range("I2").Value = 1292.8572
$ws.Range("K2").Value = 1292.8572
$ws.Range("M2").Value = -1179.8572
# Row 45
$ws.Range("H45").Value = 2878.5
$ws.Range("I45").Value = 1100
$ws.Range("J45").Value = 3363.5454
$ws.Range("K45").Value = 1100
$ws.Range("L45").Value = 3363.5454
$ws.Range("M45").Value = -723
$ws.Range("N45").Value = -4117.5454
# Row 63
$ws.Range("H63").Value = 4988.6665
$ws.Range("I63").Value = 3983
$ws.Range("K63").Value = 3983
$ws.Range("M63").Value = -3297
# Row 66
$ws.Range("H66").Value = 4988.6665
$ws.Range("I66").Value = 3983
$ws.Range("K66").Value = 19915
$ws.Range("M66").Value = -16483
# Row 75
$ws.Range("H75").Value = 23000
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
# Row 78
$ws.Range("H78").Value = 23000
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
# Row 102
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()
# Row 116
$ws.Range("H116").Value = 2672
$ws.Range("I116").Value = 1292.8572
$ws.Range("K116").Value = 1292.8572
$ws.Range("M116").Value = 1001.1428
# Row 132
$ws.Range("H132").Value = 2299.7273
$ws.Range("I132").Value = 2129.7
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 6389.099999999999
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -3859.099999999999
$ws.Range("N132").Value = -17060

# ==== Sheet: BSM ====
$ws = $wb.Worksheets.Item("BSM")

# Row 3
$ws.Range("H3").Value = 2672
$ws.Range("I3").Value = 1292.8572
$ws.Range("K3").Value = 1292.8572
$ws.Range("M3").Value = -1178.8572
# Row 99
$ws.Range("H99").Value = 2893.5833
$ws.Range("I99").Value = 2793.0908
$ws.Range("K99").Value = 2793.0908
$ws.Range("M99").Value = -1295.0908
# Row 107
$ws.Range("H107").Value = 665.7037
$ws.Range("I107").Value = 664.38464
$ws.Range("K107").Value = 664.38464
$ws.Range("M107").Value = 1255.61536

# ==== Sheet: CRP ====
$ws = $wb.Worksheets.Item("CRP")

# Row 22
$ws.Range("H22").Value = 10516.667
$ws.Range("I22").Value = 12520
$ws.Range("K22").Value = 12520
$ws.Range("M22").Value = -12170
# Row 41
$ws.Range("H41").Value = 3825.75
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
# Row 58
$ws.Range("H58").Value = 3801
$ws.Range("I58").Value = 3216.3333
$ws.Range("J58").Value = 5555
$ws.Range("K58").Value = 3216.3333
$ws.Range("L58").Value = 5555
$ws.Range("M58").Value = -3013.3333
$ws.Range("N58").Value = -5961
# Row 107
$ws.Range("H107").Value = 2317.55
$ws.Range("I107").Value = 1814.5
$ws.Range("K107").Value = 1814.5
$ws.Range("M107").Value = 105.5
# Row 132
$ws.Range("H132").Value = 2641.75
$ws.Range("I132").Value = 2641.75
$ws.Range("K132").Value = 7925.25
$ws.Range("M132").Value = -5395.25
# Row 134
$ws.Range("H134").Value = 1958.2941
$ws.Range("I134").Value = 1952.7333
$ws.Range("K134").Value = 5858.199900000001
$ws.Range("M134").Value = -3323.199900000001
# Row 136
$ws.Range("H136").Value = 3801
$ws.Range("I136").Value = 3216.3333
$ws.Range("J136").Value = 5555
$ws.Range("K136").Value = 9648.999899999999
$ws.Range("L136").Value = 16665
$ws.Range("M136").Value = -7098.999899999999
$ws.Range("N136").Value = -21765

# ==== Sheet: CUL ====
$ws = $wb.Worksheets.Item("CUL")

# Row 23
$ws.Range("H23").Value = 116.5
$ws.Range("I23").Value = 300
$ws.Range("J23").Value = 90.28570999999999
$ws.Range("K23").Value = 900
$ws.Range("L23").Value = 270.85713
$ws.Range("M23").Value = -665
$ws.Range("N23").Value = -740.85713
# Row 46
$ws.Range("H46").Value = 693.8889
$ws.Range("I46").Value = 424.25
$ws.Range("J46").Value = 909.6
$ws.Range("K46").Value = 1272.75
$ws.Range("L46").Value = 2728.8
$ws.Range("M46").Value = -1181.75
$ws.Range("N46").Value = -2910.8
# Row 56
$ws.Range("H56").Value = 8601.200000000001
$ws.Range("I56").Value = 8601.200000000001
$ws.Range("K56").Value = 8601.200000000001
$ws.Range("M56").Value = -8071.200000000001
# Row 103
$ws.Range("H103").Value = 23315.455
$ws.Range("I103").Value = 35948.855
$ws.Range("J103").Value = 1207
$ws.Range("K103").Value = 107846.565
$ws.Range("L103").Value = 3621
$ws.Range("M103").Value = -106967.565
$ws.Range("N103").Value = -5379
# Row 117
$ws.Range("H117").Value = 2308.2727
$ws.Range("I117").Value = 1819.6
$ws.Range("J117").Value = 2715.5
$ws.Range("K117").Value = 5458.799999999999
$ws.Range("L117").Value = 8146.5
$ws.Range("M117").Value = -2016.799999999999
$ws.Range("N117").Value = -15030.5
# Row 132
$ws.Range("H132").Value = 3000
$ws.Range("J132").Value = 3000
$ws.Range("L132").Value = 27000
$ws.Range("N132").Value = -32060

# ==== Sheet: GSM ====
$ws = $wb.Worksheets.Item("GSM")

# Row 12
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
# Row 70
$ws.Range("H70").Value = 4116.5
$ws.Range("I70").Value = 4116.5
$ws.Range("K70").Value = 4116.5
$ws.Range("M70").Value = -3846.5
# Row 73
$ws.Range("H73").Value = 4116.5
$ws.Range("I73").Value = 4116.5
$ws.Range("K73").Value = 4116.5
$ws.Range("M73").Value = -3180.5

# ==== Sheet: LTW ====
$ws = $wb.Worksheets.Item("LTW")

# Row 7
$ws.Range("H7").Value = 5668
$ws.Range("J7").Value = 5000
$ws.Range("L7").Value = 5000
$ws.Range("N7").Value = -5224
# Row 40
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
# Row 61
$ws.Range("H61").Value = 3314.8333
$ws.Range("I61").Value = 3097.25
$ws.Range("J61").Value = 3750
$ws.Range("K61").Value = 3097.25
$ws.Range("L61").Value = 3750
$ws.Range("M61").Value = -2895.25
$ws.Range("N61").Value = -4154
# Row 113
$ws.Range("H113").Value = 3314.8333
$ws.Range("I113").Value = 3097.25
$ws.Range("J113").Value = 3750
$ws.Range("K113").Value = 3097.25
$ws.Range("L113").Value = 3750
$ws.Range("M113").Value = -927.25
$ws.Range("N113").Value = -8090
# Row 126
$ws.Range("H126").Value = 5668
$ws.Range("J126").Value = 5000
$ws.Range("L126").Value = 15000
$ws.Range("N126").Value = -19940
# Row 132
$ws.Range("H132").Value = 3973.3
$ws.Range("I132").Value = 3859.2222
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 11577.6666
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -9047.6666
$ws.Range("N132").Value = -20060

# ==== Sheet: WVR ====
$ws = $wb.Worksheets.Item("WVR")

# Row 96
$ws.Range("H96").Value = 1690.1177
$ws.Range("I96").Value = 1656.0834
$ws.Range("J96").Value = 1771.8
$ws.Range("K96").Value = 1656.0834
$ws.Range("L96").Value = 1771.8
$ws.Range("M96").Value = -283.0834
$ws.Range("N96").Value = -4517.8
# Row 107
$ws.Range("H107").Value = 272.57144
$ws.Range("I107").Value = 272.57144
$ws.Range("K107").Value = 817.71432
$ws.Range("M107").Value = 1102.28568
# Row 113
$ws.Range("H113").Value = 655.6316
$ws.Range("I113").Value = 813.2
$ws.Range("J113").Value = 480.55554
$ws.Range("K113").Value = 2439.6
$ws.Range("L113").Value = 1441.66662
$ws.Range("M113").Value = -269.6000000000004
$ws.Range("N113").Value = -5781.66662

